# Generate Report for Handback
#
# The localization-status report is regenerated: the "Ready for handoff"
# status becomes "Handed back: in sync with en-US" everywhere it appears,
# the de-de row's handback datetime is refreshed (a fresh handback just
# happened), the now-stale "handback file is not the latest" error is
# cleared for both language rows, and the columns that hold the (now
# longer/shorter) text are resized to fit the new content.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# Widen the zh-cn / de-de status columns to fit the longer status text.
$overview.Columns.Item(5).ColumnWidth = 29.1443716684978
$overview.Columns.Item(6).ColumnWidth = 29.1443716684978

# --- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus

# The handback for zh-cn is now in sync, so the "handback file is not the
# latest" error detail is cleared, and the handback datetime is refreshed.
$zhcn.Range("J2").Value = "2016-07-26 07:53:28"
$zhcn.Range("O2").ClearContents()

$zhcn.Columns.Item(3).ColumnWidth = 29.1443716684978
$zhcn.Columns.Item(15).ColumnWidth = 12.9137198130290

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus

# A new handback just completed for de-de, so the handback datetime moves
# forward and the "not the latest" error clears the same way as zh-cn.
$dede.Range("J2").Value = "2016-07-26 07:53:43"
$dede.Range("O2").ClearContents()

$dede.Columns.Item(3).ColumnWidth = 29.1443716684978
$dede.Columns.Item(15).ColumnWidth = 12.9137198130290
